$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of player-stat data (rows 2-7), replacing the single old sample
# row. "Fecha y hora" values are written as plain text strings (quoted
# timestamps), matching the source data's text representation rather than
# native Excel date serials.
$ws.Range("A2").Value = "Solitario"
$ws.Range("B2").Value = "Annabelle"
$ws.Range("C2").Value = 190
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = "Victoria"
$ws.Range("F2").Value = "2025-11-26 15:26:59"

$ws.Range("A3").Value = "Solitario"
# B3 ("Nombre") is blank for this row.
$ws.Range("C3").Value = 288
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "Victoria"
$ws.Range("F3").Value = "2025-11-26 15:27:27"

$ws.Range("A4").Value = "Solitario"
$ws.Range("B4").Value = "Annabelle"
$ws.Range("C4").Value = 348
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = "Victoria"
$ws.Range("F4").Value = "2025-11-26 15:43:24"

$ws.Range("A5").Value = "Solitario"
$ws.Range("B5").Value = "Ulises"
$ws.Range("C5").Value = 795
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = "Derrota"
$ws.Range("F5").Value = "2025-11-26 15:43:45"

$ws.Range("A6").Value = "Multijugador"
$ws.Range("B6").Value = "Uli"
$ws.Range("C6").Value = 192
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Victoria"
$ws.Range("F6").Value = "2025-11-26 15:55:54"

$ws.Range("A7").Value = "Multijugador"
$ws.Range("B7").Value = "Pepe"
$ws.Range("C7").Value = 923
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = "Derrota"
$ws.Range("F7").Value = "2025-11-26 15:56:40"
